# Update gh-pages output (generated at 456a3b4)
# Updates the "想去人数" (want-to-go count) figures on the 展览 (Exhibition)
# and 演出 (Performance) sheets, propagating the same figures into the
# combined 全部类型 (All types) sheet, and marks a couple of shows as
# sold-out / unavailable on their respective sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 13177
$wsExpo.Range("F3").Value = 317
$wsExpo.Range("F4").Value = 644
$wsExpo.Range("F5").Value = 211
$wsExpo.Range("F6").Value = 419
$wsExpo.Range("F7").Value = 1264
$wsExpo.Range("F8").Value = 120

# ---- Sheet "演出" (performances) ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 29
$wsShow.Range("G2").Value = "不可售"

# ---- Sheet "全部类型" (all types, combined listing) ----
$wsAll = $wb.Worksheets.Item("全部类型")
# Row 2 corresponds to the 演出 sheet's row 2 entry
$wsAll.Range("F2").Value = 29
$wsAll.Range("G2").Value = "已停售"
# Rows 3-6, 9, 10, 12 correspond to the 展览 sheet's rows 2-8 entries
$wsAll.Range("F3").Value = 13177
$wsAll.Range("F4").Value = 317
$wsAll.Range("F5").Value = 644
$wsAll.Range("F6").Value = 211
$wsAll.Range("F9").Value = 419
$wsAll.Range("F10").Value = 1264
$wsAll.Range("F12").Value = 120
